$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
